$d = $word.ActiveDocument
Write-Output $d.Content.Font.Bold
